# Fruta / hortaliza, semanal
#
# A new weekly price block (3 quality grades: Especial / Primera / Segunda)
# is inserted at the top of the data table (rows 165-167), pushing every
# existing data row down by 3. The sheet dimension grows from
# A1:T183 to A1:T186 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 165; the rows below (old 165:183) shift down
# to become 168:186, retaining all of their original content untouched.
$ws.Rows("165:167").Insert()

# Row 165: Especial, new week (2023-10-12 == serial 45211)
$ws.Range("A165").Value = 2
$ws.Range("B165").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 45211
$ws.Range("E165").Value = 4
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100107
$ws.Range("H165").Value = "Otros"
$ws.Range("I165").Value = 100107002
$ws.Range("J165").Value = "Chirimoya"
$ws.Range("K165").Value = "Cultivar IV Región"
$ws.Range("L165").Value = "Especial"
$ws.Range("M165").Value = 360
$ws.Range("N165").Value = 20000
$ws.Range("O165").Value = 21000
$ws.Range("P165").Value = 20500
$ws.Range("Q165").Value = "$/bandeja 10 kilos"
$ws.Range("R165").Value = "Provincia de Limarí"
$ws.Range("S165").Value = 2050
$ws.Range("T165").Value = 10

# Row 166: Primera, same new week
$ws.Range("A166").Value = 2
$ws.Range("B166").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C166").Value = "Coquimbo"
$ws.Range("D166").Value = 45211
$ws.Range("E166").Value = 4
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100107
$ws.Range("H166").Value = "Otros"
$ws.Range("I166").Value = 100107002
$ws.Range("J166").Value = "Chirimoya"
$ws.Range("K166").Value = "Cultivar IV Región"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 400
$ws.Range("N166").Value = 17000
$ws.Range("O166").Value = 18000
$ws.Range("P166").Value = 17500
$ws.Range("Q166").Value = "$/bandeja 10 kilos"
$ws.Range("R166").Value = "Provincia de Limarí"
$ws.Range("S166").Value = 1750
$ws.Range("T166").Value = 10

# Row 167: Segunda, same new week
$ws.Range("A167").Value = 2
$ws.Range("B167").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C167").Value = "Coquimbo"
$ws.Range("D167").Value = 45211
$ws.Range("E167").Value = 4
$ws.Range("F167").Value = "Fruta"
$ws.Range("G167").Value = 100107
$ws.Range("H167").Value = "Otros"
$ws.Range("I167").Value = 100107002
$ws.Range("J167").Value = "Chirimoya"
$ws.Range("K167").Value = "Cultivar IV Región"
$ws.Range("L167").Value = "Segunda"
$ws.Range("M167").Value = 360
$ws.Range("N167").Value = 14000
$ws.Range("O167").Value = 15000
$ws.Range("P167").Value = 14500
$ws.Range("Q167").Value = "$/bandeja 10 kilos"
$ws.Range("R167").Value = "Provincia de Limarí"
$ws.Range("S167").Value = 1450
$ws.Range("T167").Value = 10

# Match the date number format already used by column D (style index 2 in
# the original file) for the freshly-inserted date cells.
$ws.Range("D165:D167").NumberFormat = $ws.Range("D168").NumberFormat
